$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.127984771727
$ws.Range("C2").Value = 0.04980699389891363
$ws.Range("D2").Value = 0.2911817816964231
$ws.Range("F2").Value = 4.45806596794165
$ws.Range("G2").Value = 0.002589902121757248
$ws.Range("K2").Value = 0.6010874286359638
$ws.Range("L2").Value = 0.3091811955591197
$ws.Range("B3").Value = 1.105637921368839
$ws.Range("C3").Value = 0.04817719131482612
$ws.Range("D3").Value = 0.2788557442090536
$ws.Range("F3").Value = 4.243910909428536
$ws.Range("G3").Value = 0.002595024083569865
$ws.Range("K3").Value = 0.5835867568609387
$ws.Range("L3").Value = 0.2989295393117857
$ws.Range("B4").Value = 1.092958261286469
$ws.Range("C4").Value = 0.04714379238185806
$ws.Range("D4").Value = 0.271212546193766
$ws.Range("F4").Value = 4.112646583209823
$ws.Range("G4").Value = 0.002598330026326999
$ws.Range("K4").Value = 0.5735605423628556
$ws.Range("L4").Value = 0.2928269107341208
$ws.Range("B5").Value = 1.088052373903707
$ws.Range("C5").Value = 0.04671431055510489
$ws.Range("D5").Value = 0.2680785989382741
$ws.Range("F5").Value = 4.05920795327873
$ws.Range("G5").Value = 0.00259971786468651
$ws.Range("K5").Value = 0.5696548602924736
$ws.Range("L5").Value = 0.2903880932124991
$ws.Range("B6").Value = 1.087253509698826
$ws.Range("C6").Value = 0.04664248614516531
$ws.Range("D6").Value = 0.2675570321820686
$ws.Range("F6").Value = 4.050337555286717
$ws.Range("G6").Value = 0.002599950772706876
$ws.Range("K6").Value = 0.5690171786308866
$ws.Range("L6").Value = 0.2899860256921727
$ws.Range("B7").Value = 1.092891042187517
$ws.Range("C7").Value = 0.04713803428897378
$ws.Range("D7").Value = 0.2711703593208483
$ws.Range("F7").Value = 4.111925684085065
$ws.Range("G7").Value = 0.002598348578574296
$ws.Range("K7").Value = 0.5735071408580978
$ws.Range("L7").Value = 0.2927938256933658
$ws.Range("B8").Value = 1.120063003219741
$ws.Range("C8").Value = 0.04925173755937351
$ws.Range("D8").Value = 0.2869470587896075
$ws.Range("F8").Value = 4.384175181173532
$ws.Range("G8").Value = 0.002591634838933902
$ws.Range("K8").Value = 0.5949034281077843
$ws.Range("L8").Value = 0.305606458545725
$ws.Range("B9").Value = 1.181647220773357
$ws.Range("C9").Value = 0.05314376850875036
$ws.Range("D9").Value = 0.317310960935842
$ws.Range("F9").Value = 4.920097657918774
$ws.Range("G9").Value = 0.002579740305245792
$ws.Range("K9").Value = 0.6426081858520831
$ws.Range("L9").Value = 0.3322661313813597
$ws.Range("B10").Value = 1.232013018401034
$ws.Range("C10").Value = 0.05585808995391872
$ws.Range("D10").Value = 0.3392994685858355
$ws.Range("F10").Value = 5.315448510380037
$ws.Range("G10").Value = 0.00257176691036552
$ws.Range("K10").Value = 0.6812212562604145
$ws.Range("L10").Value = 0.3528064831212703
$ws.Range("B11").Value = 1.256051261375774
$ws.Range("C11").Value = 0.05706334903822352
$ws.Range("D11").Value = 0.3492397714931315
$ws.Range("F11").Value = 5.495734580673457
$ws.Range("G11").Value = 0.002568303826367499
$ws.Range("K11").Value = 0.6995751408188937
$ws.Range("L11").Value = 0.3623620269669203
$ws.Range("B12").Value = 1.265316927007689
$ws.Range("C12").Value = 0.05751566594065949
$ws.Range("D12").Value = 0.3529954368600272
$ws.Range("F12").Value = 5.564073222695185
$ws.Range("G12").Value = 0.00256701588377231
$ws.Range("K12").Value = 0.7066397248354974
$ws.Range("L12").Value = 0.3660111894781721
$ws.Range("B13").Value = 1.263314140251509
$ws.Range("C13").Value = 0.05741843078622466
$ws.Range("D13").Value = 0.3521869583503872
$ws.Range("F13").Value = 5.549352157419094
$ws.Range("G13").Value = 0.002567292224428576
$ws.Range("K13").Value = 0.7051131381851974
$ws.Range("L13").Value = 0.3652239082154267
$ws.Range("B14").Value = 1.256810283150344
$ws.Range("C14").Value = 0.05710064247965896
$ws.Range("D14").Value = 0.3495489202904594
$ws.Range("F14").Value = 5.501355445220327
$ws.Range("G14").Value = 0.002568197397111671
$ws.Range("K14").Value = 0.7001540507188793
$ws.Range("L14").Value = 0.3626616293138198
$ws.Range("B15").Value = 1.252847725106164
$ws.Range("C15").Value = 0.05690546008989017
$ws.Range("D15").Value = 0.3479319510060179
$ws.Range("F15").Value = 5.471965109772725
$ws.Range("G15").Value = 0.002568754891825132
$ws.Range("K15").Value = 0.697131389350119
$ws.Range("L15").Value = 0.3610961632453922
$ws.Range("B16").Value = 1.230464809116427
$ws.Range("C16").Value = 0.05577874292265861
$ws.Range("D16").Value = 0.3386486303008667
$ws.Range("F16").Value = 5.303675659074742
$ws.Range("G16").Value = 0.002571996520331701
$ws.Range("K16").Value = 0.6800377375629125
$ws.Range("L16").Value = 0.352186285959732
$ws.Range("B17").Value = 1.217022795773801
$ws.Range("C17").Value = 0.05508009171045458
$ws.Range("D17").Value = 0.3329379621096109
$ws.Range("F17").Value = 5.200551569176781
$ws.Range("G17").Value = 0.002574027073649398
$ws.Range("K17").Value = 0.6697539969723607
$ws.Range("L17").Value = 0.3467747475786638
$ws.Range("B18").Value = 1.209397289613889
$ws.Range("C18").Value = 0.05467546014576463
$ws.Range("D18").Value = 0.329647437752314
$ws.Range("F18").Value = 5.141278330624118
$ws.Range("G18").Value = 0.002575210443737745
$ws.Range("K18").Value = 0.6639132295880188
$ws.Range("L18").Value = 0.3436820784622938
$ws.Range("B19").Value = 1.206833600075214
$ws.Range("C19").Value = 0.05453797572101138
$ws.Range("D19").Value = 0.3285322956595564
$ws.Range("F19").Value = 5.121216320003157
$ws.Range("G19").Value = 0.002575613770389768
$ws.Range("K19").Value = 0.6619483554855208
$ws.Range("L19").Value = 0.342638364872343
$ws.Range("B20").Value = 1.218442745605699
$ws.Range("C20").Value = 0.05515475154378358
$ws.Range("D20").Value = 0.3335464797053476
$ws.Range("F20").Value = 5.21152502348798
$ws.Range("G20").Value = 0.002573809319509522
$ws.Range("K20").Value = 0.6708410356600325
$ws.Range("L20").Value = 0.3473487532120316
$ws.Range("B21").Value = 1.258716195128727
$ws.Range("C21").Value = 0.05719409434223088
$ws.Range("D21").Value = 0.3503240033491295
$ws.Range("F21").Value = 5.515451350629007
$ws.Range("G21").Value = 0.002567930890687457
$ws.Range("K21").Value = 0.7016075432456716
$ws.Range("L21").Value = 0.3634133979029031
$ws.Range("B22").Value = 1.285987193060691
$ws.Range("C22").Value = 0.05850316362804975
$ws.Range("D22").Value = 0.3612396877876733
$ws.Range("F22").Value = 5.714484827672891
$ws.Range("G22").Value = 0.002564225630512257
$ws.Range("K22").Value = 0.7223823122730266
$ws.Range("L22").Value = 0.3740915513378269
$ws.Range("B23").Value = 1.271344920939612
$ws.Range("C23").Value = 0.05780661165623791
$ws.Range("D23").Value = 0.3554181491066686
$ws.Range("F23").Value = 5.608218608572713
$ws.Range("G23").Value = 0.002566190741338183
$ws.Range("K23").Value = 0.711233064652788
$ws.Range("L23").Value = 0.3683759583305601
$ws.Range("B24").Value = 1.217800467132207
$ws.Range("C24").Value = 0.05512100708867429
$ws.Range("D24").Value = 0.3332713919371599
$ws.Range("F24").Value = 5.206563881027591
$ws.Range("G24").Value = 0.002573907716483699
$ws.Range("K24").Value = 0.6703493628602928
$ws.Range("L24").Value = 0.3470891876720685
$ws.Range("B25").Value = 1.164092370877199
$ws.Range("C25").Value = 0.05211701873201946
$ws.Range("D25").Value = 0.3091550793943156
$ws.Range("F25").Value = 4.774860452769275
$ws.Range("G25").Value = 0.00258282298352518
$ws.Range("K25").Value = 0.6290813242991931
$ws.Range("L25").Value = 0.3248878771079546
